# Update Release-Notes.xlsx - Folder inventory updated on Wed Jun 11 14:47:06 UTC 2025
#
# A new folder ("Azure Landing Zone") was created after the previous run, so a
# new row is inserted at the top of the data (row 2) of the "Folder Inventory"
# sheet, pushing every existing entry down by one row. The Metadata and
# Summary sheets are then refreshed to reflect the new totals/timestamps.

$wb = $excel.ActiveWorkbook

$wsInventory = $wb.Worksheets.Item("Folder Inventory")
$wsMetadata  = $wb.Worksheets.Item("Metadata")
$wsSummary   = $wb.Worksheets.Item("Summary")

# --- Folder Inventory: insert the new row right under the header row ---
$wsInventory.Rows.Item(2).Insert()
# Excel's default Insert() copies the format of the row above (the bold
# header). The data rows use the default/no style, so clear the formatting
# that was just inherited before writing the new values.
$wsInventory.Rows.Item(2).ClearFormats()

$wsInventory.Cells.Item(2, 1).Value = "Azure Landing Zone"
$wsInventory.Cells.Item(2, 2).Value = "Azure Landing Zone"
$wsInventory.Cells.Item(2, 3).Value = "2025-06-11 20:16:49 +0530"
$wsInventory.Cells.Item(2, 4).Value = 1
$wsInventory.Cells.Item(2, 5).Value = "Root"

# --- Metadata sheet updates ---
$wsMetadata.Cells.Item(3, 2).Value = "2025-06-11 14:47:06 UTC"
$wsMetadata.Cells.Item(4, 2).Value = 69

# "Workflow Run" is stored as text ("2"), not a number - force text format
# so COM doesn't coerce the numeric-looking string back into a number, then
# drop the now-unneeded number format so the cell keeps the default style.
$wsMetadata.Cells.Item(5, 2).NumberFormat = "@"
$wsMetadata.Cells.Item(5, 2).Value = "2"
$wsMetadata.Cells.Item(5, 2).ClearFormats()

# --- Summary sheet updates ---
$wsSummary.Cells.Item(2, 2).Value = 69
$wsSummary.Cells.Item(3, 2).Value = 69
$wsSummary.Cells.Item(5, 2).Value = "2025-06-11 20:16:49 +0530"
